# Fix heat rate modeling syntax
# Updates dispatch/state-of-charge/cost output cells across several sheets
# of the "Year 0" output workbook to reflect the corrected heat-rate logic.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("D2").Value = 0
$ws.Range("S2").Value = 2.883999999999991

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76271.06239999995
$ws.Range("D2").Value = 9300.638068405267
$ws.Range("F2").Value = -1999.575668405316

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("H2").Value = 40.03636363636365
$ws.Range("L2").Value = 83.2
$ws.Range("I3").Value = 0
$ws.Range("L3").Value = 93.59999999999999
$ws.Range("M3").Value = 42.2531170288747
$ws.Range("Q3").Value = 26
$ws.Range("R3").Value = 31.2
$ws.Range("L4").Value = 29.58312417100293
$ws.Range("M4").Value = 83.2
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 72.8
$ws.Range("P4").Value = 0

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("H2").Value = 27.03636363636365
$ws.Range("L2").Value = 62.4
$ws.Range("I3").Value = 0
$ws.Range("L3").Value = 93.59999999999999
$ws.Range("M3").Value = 18.8531170288747
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 31.2
$ws.Range("L4").Value = 29.58312417100293
$ws.Range("M4").Value = 59.8
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 72.8
$ws.Range("P4").Value = 0

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("D2").Value = 13
$ws.Range("S2").Value = 7.516000000000011

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 188.6909090909091
$ws.Range("C2").Value = 168.9939393939394
$ws.Range("D2").Value = 155.8626262626262
$ws.Range("E2").Value = 142.7313131313131
$ws.Range("H2").Value = 169.236
$ws.Range("I2").Value = 189.828
$ws.Range("J2").Value = 212.994
$ws.Range("K2").Value = 259.326
$ws.Range("S2").Value = 640.4080808080807
$ws.Range("T2").Value = 608.8929292929292
$ws.Range("U2").Value = 490.7111111111111
$ws.Range("V2").Value = 392.2262626262626
$ws.Range("W2").Value = 313.4383838383839
$ws.Range("X2").Value = 260.9131313131313
$ws.Range("Y2").Value = 221.5191919191919
$ws.Range("I3").Value = 129.6
$ws.Range("J3").Value = 191.376
$ws.Range("K3").Value = 191.376
$ws.Range("L3").Value = 284.04
$ws.Range("M3").Value = 302.704585858586
$ws.Range("N3").Value = 359.332585858586
$ws.Range("O3").Value = 431.404585858586
$ws.Range("P3").Value = 454.570585858586
$ws.Range("Q3").Value = 454.570585858586
$ws.Range("B4").Value = 444.4646464646465
$ws.Range("C4").Value = 424.7676767676768
$ws.Range("D4").Value = 424.7676767676768
$ws.Range("E4").Value = 424.7676767676768
$ws.Range("F4").Value = 424.7676767676768
$ws.Range("G4").Value = 405.0707070707071
$ws.Range("H4").Value = 405.0707070707071
$ws.Range("I4").Value = 405.0707070707071
$ws.Range("J4").Value = 415.3667070707071
$ws.Range("K4").Value = 456.5507070707071
$ws.Range("L4").Value = 485.838
$ws.Range("M4").Value = 545.04
$ws.Range("N4").Value = 545.04
$ws.Range("O4").Value = 617.112
$ws.Range("P4").Value = 617.112
$ws.Range("Q4").Value = 637.704
$ws.Range("R4").Value = 648
$ws.Range("S4").Value = 648
$ws.Range("T4").Value = 516.6868686868687
$ws.Range("U4").Value = 516.6868686868687
$ws.Range("V4").Value = 516.6868686868687
$ws.Range("W4").Value = 516.6868686868687
$ws.Range("X4").Value = 516.6868686868687
$ws.Range("Y4").Value = 477.2929292929293
